$wb = $excel.ActiveWorkbook

# Column F ("想去人数") updates - same edits apply to both the "展览" and
# "全部类型" worksheets, which contain duplicated data.
$updates = @{
    2  = 163
    3  = 7359
    4  = 5739
    5  = 86
    11 = 118
    13 = 78
    14 = 657
    15 = 431
    16 = 54
    20 = 61
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
